$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '66.500.38'
Set-TextValue 'E2' '  +0.96%  '
Set-TextValue 'D3' '3.279.45'
Set-TextValue 'E3' '  +3.57%  '
Set-TextValue 'E4' '  +0.07%  '
Set-TextValue 'D5' '616.31'
Set-TextValue 'E5' '  +2.46%  '
Set-TextValue 'D6' '158.80'
Set-TextValue 'E6' '  +3.28%  '
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'D8' '3.279.94'
Set-TextValue 'E8' '  +3.53%  '
Set-TextValue 'D9' '0.546'
Set-TextValue 'E9' '  +0.63%  '
Set-TextValue 'E10' '  +4.02%  '
Set-TextValue 'E11' '  +2.16%  '
Set-TextValue 'E12' '  -2.96%  '
Set-TextValue 'D13' '0.0000273'
Set-TextValue 'E13' '  +3.52%  '
Set-TextValue 'D14' '39.25'
Set-TextValue 'E14' '  +3.20%  '
Set-TextValue 'D15' '3.818.26'
Set-TextValue 'E15' '  +3.65%  '
Set-TextValue 'D16' '66.568.93'
Set-TextValue 'E16' '  +0.90%  '
Set-TextValue 'E17' '  +1.48%  '
Set-TextValue 'D18' '3.281.44'
Set-TextValue 'E18' '  +3.49%  '
Set-TextValue 'E19' '  +1.55%  '
Set-TextValue 'D20' '507.00'
Set-TextValue 'D21' '15.60'
Set-TextValue 'E21' '  +2.27%  '
Set-TextValue 'E22' '  +4.58%  '
Set-TextValue 'D23' '8.20'
Set-TextValue 'E23' '  +3.19%  '
Set-TextValue 'D24' '14.75'
Set-TextValue 'E24' '  +0.23%  '
Set-TextValue 'D25' '86.96'
Set-TextValue 'E25' '  +3.28%  '
Set-TextValue 'E26' '  +0.07%  '
Set-TextValue 'E27' '  +2.87%  '
Set-TextValue 'D28' '9.31'
Set-TextValue 'E28' '  +2.21%  '
Set-TextValue 'D29' '2.43'
Set-TextValue 'E29' '  +2.55%  '
Set-TextValue 'D30' '0.129'
Set-TextValue 'E30' '  +46.91%  '
Set-TextValue 'E31' '  -1.73%  '
Set-TextValue 'E32' '  -2.66%  '
Set-TextValue 'D33' '28.16'
Set-TextValue 'E33' '  +1.22%  '
Set-TextValue 'D34' '1.00'
Set-TextValue 'E34' '  -0.16%  '
Set-TextValue 'E35' '  -2.03%  '
Set-TextValue 'D36' '6.51'
Set-TextValue 'E36' '  +0.80%  '
Set-TextValue 'D37' '3.44'
Set-TextValue 'E37' '  +21.22%  '
Set-TextValue 'B38' 'OKB'
Set-TextValue 'C38' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D38' '55.68'
Set-TextValue 'E38' '  +0.79%  '
Set-TextValue 'B39' 'PEPE'
Set-TextValue 'C39' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D39' '0.0₃0793'
Set-TextValue 'E39' '  +17.53%  '
Set-TextValue 'D40' '496.69'
Set-TextValue 'E40' '  -1.67%  '
Set-TextValue 'E41' '  +2.83%  '
Set-TextValue 'E42' '  +2.00%  '
Set-TextValue 'E43' '  +1.55%  '
Set-TextValue 'D44' '2.61'
Set-TextValue 'E44' '  +7.67%  '
Set-TextValue 'D45' '0.297'
Set-TextValue 'E45' '  +0.42%  '
Set-TextValue 'D46' '3.021.30'
Set-TextValue 'E46' '  +7.18%  '
Set-TextValue 'D47' '29.28'
Set-TextValue 'E47' '  +5.56%  '
Set-TextValue 'E48' '  +6.49%  '
Set-TextValue 'E49' '  +3.06%  '
Set-TextValue 'E50' '  -0.01%  '
Set-TextValue 'B51' 'Monero'
Set-TextValue 'C51' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D51' '121.61'
Set-TextValue 'E51' '  -0.02%  '
